# Auto-generated Excel COM-interop script
# Applies cell value updates per the target diff (market price refresh).
$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 622.38464
$ws.Range("I2").Value = 622.38464
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 622.38464
$ws.Range("L2").Value = 0
$ws.Range("M2").Value = -509.38464
$ws.Range("N2").ClearContents()
$ws.Range("H32").Value = 5888.067
$ws.Range("J32").Value = 5768.7144
$ws.Range("L32").Value = 5768.7144
$ws.Range("N32").Value = -6420.7144
$ws.Range("H112").Value = 1442.1143
$ws.Range("J112").Value = 1493.2307
$ws.Range("L112").Value = 4479.6921
$ws.Range("N112").Value = -6695.6921
$ws.Range("H137").Value = 7729.6387
$ws.Range("I137").Value = 9052.607
$ws.Range("K137").Value = 27157.821
$ws.Range("M137").Value = -24607.821
$ws.Range("H138").Value = 3341.6123
$ws.Range("J138").Value = 4390
$ws.Range("L138").Value = 13170
$ws.Range("N138").Value = -23450
$ws.Range("H141").Value = 4868.1724
$ws.Range("I141").Value = 4065.3809
$ws.Range("K141").Value = 12196.1427
$ws.Range("M141").Value = -7016.1427

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 4567.6216
$ws.Range("I61").Value = 3821.1538
$ws.Range("K61").Value = 3821.1538
$ws.Range("M61").Value = -3609.1538
$ws.Range("H80").Value = 85880
$ws.Range("J80").Value = 85880
$ws.Range("L80").Value = 85880
$ws.Range("N80").Value = -87876
$ws.Range("H83").Value = 85880
$ws.Range("J83").Value = 85880
$ws.Range("L83").Value = 257640
$ws.Range("N83").Value = -267624
$ws.Range("H88").Value = 3791.3333
$ws.Range("I88").Value = 2639.6
$ws.Range("J88").Value = 4614
$ws.Range("K88").Value = 2639.6
$ws.Range("L88").Value = 4614
$ws.Range("M88").Value = -2233.6
$ws.Range("N88").Value = -5426
$ws.Range("H91").Value = 3791.3333
$ws.Range("I91").Value = 2639.6
$ws.Range("J91").Value = 4614
$ws.Range("K91").Value = 2639.6
$ws.Range("L91").Value = 4614
$ws.Range("M91").Value = -1235.6
$ws.Range("N91").Value = -7422
$ws.Range("H122").Value = 643094.5600000001
$ws.Range("I122").Value = 6469.857
$ws.Range("K122").Value = 19409.571
$ws.Range("M122").Value = -16959.571
$ws.Range("H136").Value = 4567.6216
$ws.Range("I136").Value = 3821.1538
$ws.Range("K136").Value = 11463.4614
$ws.Range("M136").Value = -8913.4614

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H32").Value = 52895
$ws.Range("J32").Value = 52895
$ws.Range("L32").Value = 52895
$ws.Range("N32").Value = -53663
$ws.Range("H82").Value = 63609.75
$ws.Range("I82").Value = 29440
$ws.Range("J82").Value = 74999.664
$ws.Range("K82").Value = 29440
$ws.Range("L82").Value = 74999.664
$ws.Range("M82").Value = -29057
$ws.Range("N82").Value = -75765.664
$ws.Range("H85").Value = 63609.75
$ws.Range("I85").Value = 29440
$ws.Range("J85").Value = 74999.664
$ws.Range("K85").Value = 29440
$ws.Range("L85").Value = 74999.664
$ws.Range("M85").Value = -28114
$ws.Range("N85").Value = -77651.664
$ws.Range("H86").Value = 8880.727999999999
$ws.Range("I86").Value = 14560.2
$ws.Range("K86").Value = 14560.2
$ws.Range("M86").Value = -13437.2
$ws.Range("H89").Value = 8880.727999999999
$ws.Range("I89").Value = 14560.2
$ws.Range("K89").Value = 72801
$ws.Range("M89").Value = -67185
$ws.Range("H99").Value = 14149.883
$ws.Range("I99").Value = 14581.125
$ws.Range("J99").Value = 7250
$ws.Range("K99").Value = 14581.125
$ws.Range("L99").Value = 7250
$ws.Range("M99").Value = -13083.125
$ws.Range("N99").Value = -10246
$ws.Range("H103").Value = 23661.666
$ws.Range("J103").Value = 23661.666
$ws.Range("L103").Value = 23661.666
$ws.Range("N103").Value = -26005.666
$ws.Range("H134").Value = 1916.4474
$ws.Range("I134").Value = 1268.1818
$ws.Range("J134").Value = 6195
$ws.Range("K134").Value = 3804.5454
$ws.Range("L134").Value = 18585
$ws.Range("M134").Value = -1269.5454
$ws.Range("N134").Value = -23655

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H19").Value = 269.8125
$ws.Range("I19").Value = 165.63637
$ws.Range("J19").Value = 499
$ws.Range("K19").Value = 165.63637
$ws.Range("L19").Value = 499
$ws.Range("M19").Value = 4.363630000000001
$ws.Range("N19").Value = -839
$ws.Range("H22").Value = 1191.2307
$ws.Range("I22").Value = 619
$ws.Range("J22").Value = 1445.5555
$ws.Range("K22").Value = 619
$ws.Range("L22").Value = 1445.5555
$ws.Range("M22").Value = -269
$ws.Range("N22").Value = -2145.5555
$ws.Range("H24").Value = 269.8125
$ws.Range("I24").Value = 165.63637
$ws.Range("J24").Value = 499
$ws.Range("K24").Value = 165.63637
$ws.Range("L24").Value = 499
$ws.Range("M24").Value = 4.363630000000001
$ws.Range("N24").Value = -839

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H47").Value = 478.6
$ws.Range("I47").Value = 198.33333
$ws.Range("J47").Value = 899
$ws.Range("K47").Value = 594.99999
$ws.Range("L47").Value = 2697
$ws.Range("M47").Value = -163.99999
$ws.Range("N47").Value = -3559
$ws.Range("H50").Value = 803.2308
$ws.Range("I50").Value = 403.81818
$ws.Range("J50").Value = 3000
$ws.Range("K50").Value = 1211.45454
$ws.Range("L50").Value = 9000
$ws.Range("M50").Value = -730.45454
$ws.Range("N50").Value = -9962
$ws.Range("H53").Value = 803.2308
$ws.Range("I53").Value = 403.81818
$ws.Range("J53").Value = 3000
$ws.Range("K53").Value = 1211.45454
$ws.Range("L53").Value = 9000
$ws.Range("M53").Value = -730.45454
$ws.Range("N53").Value = -9962
$ws.Range("H122").Value = 2273.2856
$ws.Range("J122").Value = 3388.4666
$ws.Range("L122").Value = 30496.1994
$ws.Range("N122").Value = -35396.1994
$ws.Range("H134").Value = 3992.6428
$ws.Range("I134").Value = 2377.4443
$ws.Range("J134").Value = 6900
$ws.Range("K134").Value = 7132.3329
$ws.Range("L134").Value = 20700
$ws.Range("M134").Value = -2062.3329
$ws.Range("N134").Value = -30840

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H26").Value = 24999
$ws.Range("J26").Value = 24999
$ws.Range("L26").Value = 24999
$ws.Range("N26").Value = -25559
$ws.Range("H50").Value = 24999
$ws.Range("J50").Value = 24999
$ws.Range("L50").Value = 24999
$ws.Range("N50").Value = -25995
$ws.Range("H102").Value = 7215.5
$ws.Range("I102").Value = 9245.412
$ws.Range("K102").Value = 9245.412
$ws.Range("M102").Value = -7623.412
$ws.Range("H126").Value = 13031.556
$ws.Range("I126").Value = 18322.1
$ws.Range("J126").Value = 9919.471
$ws.Range("K126").Value = 54966.3
$ws.Range("L126").Value = 29758.413
$ws.Range("M126").Value = -52496.3
$ws.Range("N126").Value = -34698.413

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 28387.834
$ws.Range("I7").Value = 49720.555
$ws.Range("K7").Value = 49720.555
$ws.Range("M7").Value = -49608.555
$ws.Range("H22").Value = 3694.6924
$ws.Range("J22").Value = 3499.5
$ws.Range("L22").Value = 3499.5
$ws.Range("N22").Value = -4089.5
$ws.Range("H27").Value = 3694.6924
$ws.Range("J27").Value = 3499.5
$ws.Range("L27").Value = 3499.5
$ws.Range("N27").Value = -3713.5
$ws.Range("H40").Value = 18648.592
$ws.Range("I40").Value = 29677.4
$ws.Range("J40").Value = 4862.5835
$ws.Range("K40").Value = 29677.4
$ws.Range("L40").Value = 4862.5835
$ws.Range("M40").Value = -29541.4
$ws.Range("N40").Value = -5134.5835
$ws.Range("H46").Value = 3972.0908
$ws.Range("I46").Value = 939.8
$ws.Range("J46").Value = 6499
$ws.Range("K46").Value = 939.8
$ws.Range("L46").Value = 6499
$ws.Range("M46").Value = -751.8
$ws.Range("N46").Value = -6875
$ws.Range("H55").Value = 835.1429000000001
$ws.Range("I55").Value = 835.1429000000001
$ws.Range("K55").Value = 835.1429000000001
$ws.Range("M55").Value = -662.1429000000001
$ws.Range("H56").Value = 33437
$ws.Range("I56").Value = 33437
$ws.Range("K56").Value = 33437
$ws.Range("M56").Value = -32746
$ws.Range("H59").Value = 44940
$ws.Range("J59").Value = 44940
$ws.Range("L59").Value = 44940
$ws.Range("N59").Value = -46248
$ws.Range("H93").Value = 5672.7617
$ws.Range("I93").Value = 5672.7617
$ws.Range("K93").Value = 5672.7617
$ws.Range("M93").Value = -4424.7617
$ws.Range("H126").Value = 28387.834
$ws.Range("I126").Value = 49720.555
$ws.Range("K126").Value = 149161.665
$ws.Range("M126").Value = -146691.665
$ws.Range("H136").Value = 10843.883
$ws.Range("I136").Value = 27284.666
$ws.Range("J136").Value = 7320.857
$ws.Range("K136").Value = 81853.99800000001
$ws.Range("L136").Value = 21962.571
$ws.Range("M136").Value = -79303.99800000001
$ws.Range("N136").Value = -27062.571

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H61").Value = 27473.75
$ws.Range("I61").Value = 27473.75
$ws.Range("K61").Value = 27473.75
$ws.Range("M61").Value = -27181.75
$ws.Range("H113").Value = 2742.875
$ws.Range("I113").Value = 1206.8422
$ws.Range("K113").Value = 3620.5266
$ws.Range("M113").Value = -1450.5266
$ws.Range("H136").Value = 2514.389
$ws.Range("I136").Value = 1808.4615
$ws.Range("K136").Value = 5425.3845
$ws.Range("M136").Value = -2875.3845

Write-Host "Applied 242 cell updates across 8 sheets."
